$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.255.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.907.35"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.08%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5251"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3793"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.69%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07271"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.47%  "
$ws.Range("E10").Value = "  +3.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9004"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08141"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +8.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.905.01"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "95.42"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.297"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008640"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.51"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.313.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.070"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.142.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.459"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.315"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +11.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "146.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.752"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.993"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.814"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09234"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8090"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05064"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.245"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.005"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.331"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.585"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5739"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01994"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.080"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.73%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "119.60"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.638"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.980"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1517"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4862"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.001"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.627"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.80%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.87"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.56%  "
